# Update Correspond Handoff/Handback datetimes for the "18d96c8f..." row
# on both the zh-cn and de-de sheets, as part of regenerating the
# Handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 06:24:31"
$wsZhCn.Range("H3").Value = "2016-03-22 06:25:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 06:24:38"
$wsDeDe.Range("H3").Value = "2016-03-22 06:25:45"
